$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("展览")
$ws.Range("F7").Value = 2590
$ws.Range("F8").Value = 37
$ws.Range("F10").Value = 915
$ws.Range("F12").Value = 920
$ws.Range("F13").Value = 1154
$ws.Range("F17").Value = 737
$ws.Range("F18").Value = 785
$ws.Range("F19").Value = 212
$ws.Range("F20").Value = 502
$ws.Range("F21").Value = 1126
$ws.Range("F23").Value = 620
$ws.Range("F26").Value = 308
$ws.Range("F29").Value = 517
$ws.Range("F30").Value = 4550
$ws.Range("F31").Value = 4550
$ws.Range("F36").Value = 164
$ws.Range("F39").Value = 64
$ws.Range("F42").Value = 88
$ws.Range("F47").Value = 133
$ws.Range("F48").Value = 110

$ws = $wb.Worksheets.Item("演出")
$ws.Range("F9").Value = 30
$ws.Range("F14").Value = 28
$ws.Range("F17").Value = 199
$ws.Range("F20").Value = 1

$ws = $wb.Worksheets.Item("本地生活")
$ws.Range("F3").Value = 731

$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F6").Value = 2590
$ws.Range("F7").Value = 37
$ws.Range("F9").Value = 915
$ws.Range("F11").Value = 920
$ws.Range("F12").Value = 1154
$ws.Range("F16").Value = 737
$ws.Range("F18").Value = 785
$ws.Range("F19").Value = 212
$ws.Range("F20").Value = 502
$ws.Range("F21").Value = 1126
$ws.Range("F25").Value = 30
$ws.Range("F26").Value = 620
$ws.Range("F28").Value = 308
$ws.Range("F30").Value = 517
$ws.Range("F31").Value = 4550
$ws.Range("F36").Value = 164
$ws.Range("F40").Value = 28
$ws.Range("F41").Value = 28
$ws.Range("F43").Value = 88
$ws.Range("F49").Value = 110
